$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 24.78000000000043
$ws.Range("H2").Value = 0.002927170793914469
$ws.Range("I2").Value = 0.002927170793914469
$ws.Range("L2").Value = 38.7019213866825
$ws.Range("M2").Value = "[12.590388441705045, 64.81345433165995]"
$ws.Range("N2").Value = 0.004569479995429004
$ws.Range("O2").Value = 0.004569479995429004
$ws.Range("P2").Value = 1.691868716347656
$ws.Range("Q2").Value = "[0.798763297309117, 2.5849741353861955]"
$ws.Range("R2").Value = 0.0004120771717011085
$ws.Range("S2").Value = 0.0004120771717011085
$ws.Range("T2").Value = 56.82274686123396
$ws.Range("U2").Value = "[41.469366372470034, 72.17612734999788]"
$ws.Range("V2").Value = [double]"2.171652635496457e-09"
$ws.Range("W2").Value = [double]"2.171652635496457e-09"
$ws.Range("X2").Value = 18.10750750750783
$ws.Range("Y2").Value = 14.58522522522548
$ws.Range("Z2").Value = 21.62978978979017

$ws.Range("B3").Value = 0
$ws.Range("F3").Value = 24.78000000000043
$ws.Range("H3").Value = 0.02032794810773719
$ws.Range("I3").Value = 0.02032794810773719
$ws.Range("L3").Value = 40.26889355157563
$ws.Range("M3").Value = "[2.4970651885962383, 78.04072191455502]"
$ws.Range("N3").Value = 0.03719675012642609
$ws.Range("O3").Value = 0.03719675012642609
$ws.Range("P3").Value = 1.062921238151501
$ws.Range("Q3").Value = "[-0.20755266780473214, 2.333395144107735]"
$ws.Range("R3").Value = 0.09889962121933005
$ws.Range("S3").Value = 0.09889962121933005
$ws.Range("T3").Value = 78.77305579097774
$ws.Range("U3").Value = "[59.025608608510765, 98.5205029734447]"
$ws.Range("V3").Value = [double]"3.072868626219361e-10"
$ws.Range("W3").Value = [double]"3.072868626219361e-10"
$ws.Range("X3").Value = 20.58798798798835
$ws.Range("Y3").Value = 15.57741741741769
$ws.Range("Z3").Value = 25.59855855855901

$ws.Range("F4").Value = 24.78000000000043
$ws.Range("H4").Value = 0.000552032042940187
$ws.Range("I4").Value = 0.000552032042940187
$ws.Range("L4").Value = 40.47540408549336
$ws.Range("M4").Value = "[14.993140263636107, 65.95766790735061]"
$ws.Range("N4").Value = 0.002526144366455707
$ws.Range("O4").Value = 0.002526144366455707
$ws.Range("P4").Value = 1.050342288587578
$ws.Range("Q4").Value = "[0.3962369112635775, 1.7044476659115793]"
$ws.Range("R4").Value = 0.002287978266036106
$ws.Range("S4").Value = 0.002287978266036106
$ws.Range("T4").Value = 59.55796152866405
$ws.Range("U4").Value = "[45.86021082775358, 73.25571222957453]"
$ws.Range("V4").Value = [double]"2.798139497883767e-11"
$ws.Range("W4").Value = [double]"2.798139497883767e-11"
$ws.Range("X4").Value = 20.63759759759796
$ws.Range("Y4").Value = 18.05789789789821
$ws.Range("Z4").Value = 23.2172972972977

$ws.Range("F5").Value = 24.78000000000043
$ws.Range("H5").Value = 0.00481580088057465
$ws.Range("I5").Value = 0.00481580088057465
$ws.Range("L5").Value = 42.48786522296487
$ws.Range("M5").Value = "[10.683753713985723, 74.29197673194402]"
$ws.Range("N5").Value = 0.00997164277332141
$ws.Range("O5").Value = 0.00997164277332141
$ws.Range("P5").Value = 1.062921238151501
$ws.Range("Q5").Value = "[0.2201316173686534, 1.9057108589343494]"
$ws.Range("R5").Value = 0.01460119155763562
$ws.Range("S5").Value = 0.01460119155763562
$ws.Range("T5").Value = 60.33201655924363
$ws.Range("U5").Value = "[42.83518103221765, 77.82885208626962]"
$ws.Range("V5").Value = [double]"1.229702939475885e-08"
$ws.Range("W5").Value = [double]"1.229702939475885e-08"
$ws.Range("X5").Value = 20.58798798798835
$ws.Range("Y5").Value = 17.26414414414445
$ws.Range("Z5").Value = 23.91183183183225

$ws.Range("B6").Value = 0
$ws.Range("F6").Value = 25.40000000000053
$ws.Range("H6").Value = 0.3182479453704019
$ws.Range("I6").Value = 0.3182479453704019
$ws.Range("L6").Value = 21.82979148571616
$ws.Range("M6").Value = "[-13.710402843472433, 57.36998581490475]"
$ws.Range("N6").Value = 0.2224599450487097
$ws.Range("O6").Value = 0.2224599450487097
$ws.Range("P6").Value = -0.0503157982556921
$ws.Range("Q6").Value = "[-3.163605815326658, 3.062974218815274]"
$ws.Range("R6").Value = 0.9741764665220474
$ws.Range("S6").Value = 0.9741764665220474
$ws.Range("T6").Value = 65.46664107622792
$ws.Range("U6").Value = "[45.69510862768205, 85.23817352477377]"
$ws.Range("V6").Value = [double]"3.160205430674523e-08"
$ws.Range("W6").Value = [double]"3.160205430674523e-08"
$ws.Range("X6").Value = 0.2034034034034065
$ws.Range("Y6").Value = -12.38218218218244
$ws.Range("Z6").Value = 12.78898898898925

$ws.Range("B7").Value = 0
$ws.Range("F7").Value = 25.40000000000053
$ws.Range("H7").Value = 0.001737056516448621
$ws.Range("I7").Value = 0.001737056516448621
$ws.Range("L7").Value = 44.72107626574949
$ws.Range("M7").Value = "[14.670378177731862, 74.77177435376711]"
$ws.Range("N7").Value = 0.00442123268499417
$ws.Range("O7").Value = 0.00442123268499417
$ws.Range("P7").Value = 0.3710790121357324
$ws.Range("Q7").Value = "[-0.3585000625718058, 1.1006580868432705]"
$ws.Range("R7").Value = 0.3111150433831364
$ws.Range("S7").Value = 0.3111150433831364
$ws.Range("T7").Value = 63.57238131135789
$ws.Range("U7").Value = "[47.42311436879508, 79.72164825392069]"
$ws.Range("V7").Value = [double]"4.379490103900707e-10"
$ws.Range("W7").Value = [double]"4.379490103900707e-10"
$ws.Range("X7").Value = 23.8998998999004
$ws.Range("Y7").Value = 20.95055055055099
$ws.Range("Z7").Value = 26.8492492492498

$ws.Range("F8").Value = 25.40000000000053
$ws.Range("H8").Value = 0.001584296139746089
$ws.Range("I8").Value = 0.001584296139746089
$ws.Range("L8").Value = 37.09850568070408
$ws.Range("M8").Value = "[12.66976025500773, 61.52725110640043]"
$ws.Range("N8").Value = 0.003736870202803999
$ws.Range("O8").Value = 0.003736870202803999
$ws.Range("P8").Value = -0.2767368904063083
$ws.Range("Q8").Value = "[-1.0188949146777704, 0.4654211338651537]"
$ws.Range("R8").Value = 0.4565492699048339
$ws.Range("S8").Value = 0.4565492699048339
$ws.Range("T8").Value = 57.44017915224416
$ws.Range("U8").Value = "[43.96211480603064, 70.91824349845767]"
$ws.Range("V8").Value = [double]"4.953126797602181e-11"
$ws.Range("W8").Value = [double]"4.953126797602181e-11"
$ws.Range("X8").Value = 1.118718718718743
$ws.Range("Y8").Value = -1.881481481481518
$ws.Range("Z8").Value = 4.118918918919004

$ws.Range("B9").Value = 0
$ws.Range("F9").Value = 25.40000000000053
$ws.Range("H9").Value = 0.02641713325594897
$ws.Range("I9").Value = 0.02641713325594897
$ws.Range("L9").Value = 37.02546424368096
$ws.Range("M9").Value = "[0.9018907188376915, 73.14903776852422]"
$ws.Range("N9").Value = 0.04477341738075657
$ws.Range("O9").Value = 0.04477341738075657
$ws.Range("P9").Value = -0.69184222601577
$ws.Range("Q9").Value = "[-2.283079345852041, 0.8993948938205012]"
$ws.Range("R9").Value = 0.3858451906337668
$ws.Range("S9").Value = 0.3858451906337668
$ws.Range("T9").Value = 63.20646533241076
$ws.Range("U9").Value = "[44.24265404103944, 82.17027662378209]"
$ws.Range("V9").Value = [double]"2.718295433190576e-08"
$ws.Range("W9").Value = [double]"2.718295433190576e-08"
$ws.Range("X9").Value = 2.796796796796851
$ws.Range("Y9").Value = -3.635835835835914
$ws.Range("Z9").Value = 9.229429429429617

$ws.Range("F10").Value = 25.40000000000053
$ws.Range("H10").Value = [double]"7.73162001090899e-05"
$ws.Range("I10").Value = [double]"7.73162001090899e-05"
$ws.Range("L10").Value = 50.52648508248384
$ws.Range("M10").Value = "[21.299184318024842, 79.75378584694283]"
$ws.Range("N10").Value = 0.00111957578250399
$ws.Range("O10").Value = 0.00111957578250399
$ws.Range("P10").Value = -0.9937370155499243
$ws.Range("Q10").Value = "[-1.5597897459264631, -0.4276842851733855]"
$ws.Range("R10").Value = 0.0009549273310984763
$ws.Range("S10").Value = 0.0009549273310984763
$ws.Range("T10").Value = 68.76435248725426
$ws.Range("U10").Value = "[53.568162619900036, 83.96054235460848]"
$ws.Range("V10").Value = [double]"8.751888103120109e-12"
$ws.Range("W10").Value = [double]"8.751888103120109e-12"
$ws.Range("X10").Value = 4.017217217217301
$ws.Range("Y10").Value = 1.728928928928966
$ws.Range("Z10").Value = 6.305505505505635
